$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added to the dataset. In the source data
# (ordered newest-first) this record belongs at row 154, which pushes every
# existing record from row 154 down by one row (the former row 247 record
# ends up at row 248, growing the sheet from 247 to 248 data rows).
$ws.Rows(154).Insert()

$ws.Range("A154").Value = 6
$ws.Range("B154").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C154").Value = "Metropolitana"
$ws.Range("D154").Value = 44824
$ws.Range("E154").Value = 13
$ws.Range("F154").Value = 100112001
$ws.Range("G154").Value = "Berenjena"
$ws.Range("H154").Value = "Sin especificar"
$ws.Range("I154").Value = "Primera"
$ws.Range("J154").Value = 180
$ws.Range("K154").Value = 12000
$ws.Range("L154").Value = 12000
$ws.Range("M154").Value = 12000
$ws.Range("N154").Value = "$/caja 50 unidades"
$ws.Range("O154").Value = "Provincia de Limarí"
$ws.Range("P154").Value = 240
$ws.Range("Q154").Value = 50
$ws.Range("R154").Value = "Hortaliza"
